$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the row for student 1252675 (row 10) - shifts everything below up by one.
$ws.Rows("10").Delete()

# Clear the leftover sort-state memory from the prior sort operation.
$ws.Sort.SortFields.Clear()

# 2) Highlight everyone who replied "Yes" (e-mailed OK with moving the dates)
#    with a green fill. Rows whose answer changed from "?" to "Yes" get the
#    text updated too (leading "'" forces text/quote-prefix, matching the
#    original cells' formatting intent). B5 (0938323) was already "Yes".
$green = 6750054  # RGB(102,255,102) -> FF66FF66

$yesRows = @(4,5,6,7,10,11,12,13,14,15,16)
foreach ($r in $yesRows) {
    $ws.Range("B$r").Value2 = "'Yes"
    $ws.Range("B$r").Interior.Color = $green
}

# 3) Selection left where the user ended up after deleting the row (rows 20:23 full-row select)
$ws.Range("A20:XFD23").Select()

# 4) Stale AutoFilter defined name left behind (hidden, sheet-scoped, pointing at #REF!)
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!#REF!", $false)
$nm.Visible = $false
